# Complete test for Systematic_Risk macro
#
# The original row 27 ("Systematic Risk1" / "Test Systematic Risk" /
# "Systematic_Risk_test1") is replaced by a pair of rows appended at the
# bottom of the table: one for scale=252 (reusing the original test name)
# and a new one for scale=1. Every row that used to follow row 27 shifts
# up by one to fill the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Systematic Risk1" row; rows below shift up by one
# (old row 28 "Bull/Bear beta" becomes the new row 27, etc.)
$ws.Rows("27:27").Delete() | Out-Null

# Append the two "Systematic Risk" test rows at the end of the table
$ws.Range("A58").Value = "Systematic Risk1"
$ws.Range("B58").Value = "Test Systematic Risk with scale=252"
$ws.Range("C58").Value = "Systematic_Risk_test1"

$ws.Range("A59").Value = "Systematic Risk2"
$ws.Range("B59").Value = "Test Systematic Risk with scale=1"
$ws.Range("C59").Value = "Systematic_Risk_test2"

# Update the window's scroll position / selection to match the edited sheet
$win = $excel.ActiveWindow
$win.ScrollRow = 39
$win.ScrollColumn = 1
$ws.Range("G61").Select() | Out-Null
